$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A3").Value = 45991
$ws.Range("B3").Value = 48994

$ws.Range("A3").NumberFormat = $ws.Range("A2").NumberFormat
$ws.Range("B3").NumberFormat = $ws.Range("B2").NumberFormat
